$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.770.40"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.636.36"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.29"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.258"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0636"
$ws.Range("E9").Value = "  -1.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("E10").Value = "  -4.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("D12").Value = "1.650.37"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").Value = "1.866.58"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.554"
$ws.Range("E15").Value = "  -1.60%  "

$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.80"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").Value = "25.792.91"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.82"
$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.96"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.53"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0492"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.24"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  +1.03%  "

$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.549"
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.118.12"
$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.52"
$ws.Range("E39").Value = "  -1.11%  "

$ws.Range("E40").Value = "  -0.61%  "

$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.67"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.799"
$ws.Range("E44").Value = "  -0.57%  "

$ws.Range("D45").Value = "1.769.33"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.24"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.418"
$ws.Range("E48").Value = "  -2.27%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.64"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  +3.19%  "
